$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 5621813206
$ws.Range("B2").Value = "Saad"
$ws.Range("C2").Value = "ar"
$ws.Range("D2").Value = "'False"
$ws.Range("E2").Value = "UTC"
$ws.Range("F2").Value = "2025-11-11 11:05:59"
